# Auto-generated Excel COM-interop edit script
# Updates raw market-price snapshot values (columns H:N) on several leve rows
# across the ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR sheets, matching the
# scheduled market-data refresh described in the commit message.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 113
$ws.Range("H113").Value = 2968.5715
$ws.Range("I113").Value = 2523.077
$ws.Range("K113").Value = 2523.077
$ws.Range("M113").Value = 730.9229999999998
# Row 116
$ws.Range("H116").Value = 4666.6665
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 4666.6665
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = ""
$ws.Range("M116").Value = 4666.6665
$ws.Range("N116").Value = -11550.6665
# Row 132
$ws.Range("H132").Value = 5104205
$ws.Range("I132").Value = 6495748.5
$ws.Range("J132").Value = 1880
$ws.Range("K132").Value = 19487245.5
$ws.Range("L132").Value = 5640
$ws.Range("M132").Value = -19484715.5
$ws.Range("N132").Value = -10700
# Row 138
$ws.Range("H138").Value = 1878.5714
$ws.Range("I138").Value = 1657.8948
$ws.Range("J138").Value = 1960.7843
$ws.Range("K138").Value = 4973.6844
$ws.Range("L138").Value = 5882.3529
$ws.Range("M138").Value = 166.3155999999999
$ws.Range("N138").Value = -16162.3529

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 43479096
$ws.Range("I2").Value = 125000540
$ws.Range("J2").Value = 991.86664
$ws.Range("K2").Value = 125000540
$ws.Range("L2").Value = 991.86664
$ws.Range("M2").Value = -125000427
$ws.Range("N2").Value = -1217.86664
# Row 5
$ws.Range("H5").Value = 27777888
$ws.Range("I5").Value = 41666732
$ws.Range("J5").Value = 199.5
$ws.Range("K5").Value = 41666732
$ws.Range("L5").Value = 199.5
$ws.Range("M5").Value = -41666620
$ws.Range("N5").Value = -423.5
# Row 21
$ws.Range("H21").Value = 500
$ws.Range("I21").Value = 500
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 500
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = ""
$ws.Range("N21").Value = -126
# Row 45
$ws.Range("H45").Value = 47620450
$ws.Range("I45").Value = 55556940
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 55556940
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -55556563
$ws.Range("N45").Value = -2254
# Row 61
$ws.Range("H61").Value = 2312.4707
$ws.Range("I61").Value = 1679.1111
$ws.Range("J61").Value = 3025
$ws.Range("K61").Value = 1679.1111
$ws.Range("L61").Value = 3025
$ws.Range("M61").Value = -1467.1111
$ws.Range("N61").Value = -3449
# Row 116
$ws.Range("H116").Value = 43479096
$ws.Range("I116").Value = 125000540
$ws.Range("J116").Value = 991.86664
$ws.Range("K116").Value = 125000540
$ws.Range("L116").Value = 991.86664
$ws.Range("M116").Value = -124998246
$ws.Range("N116").Value = -5579.86664
# Row 136
$ws.Range("H136").Value = 2312.4707
$ws.Range("I136").Value = 1679.1111
$ws.Range("J136").Value = 3025
$ws.Range("K136").Value = 5037.3333
$ws.Range("L136").Value = 9075
$ws.Range("M136").Value = -2487.3333
$ws.Range("N136").Value = -14175

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 43479096
$ws.Range("I3").Value = 125000540
$ws.Range("J3").Value = 991.86664
$ws.Range("K3").Value = 125000540
$ws.Range("L3").Value = 991.86664
$ws.Range("M3").Value = -125000426
$ws.Range("N3").Value = -1219.86664
# Row 4
$ws.Range("H4").Value = 27777888
$ws.Range("I4").Value = 41666732
$ws.Range("J4").Value = 199.5
$ws.Range("K4").Value = 41666732
$ws.Range("L4").Value = 199.5
$ws.Range("M4").Value = -41666617
$ws.Range("N4").Value = -429.5
# Row 25
$ws.Range("H25").Value = 1207
$ws.Range("I25").Value = 1207
$ws.Range("K25").Value = 1207
$ws.Range("M25").Value = -972

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 870.6957
$ws.Range("I16").Value = 847.94446
$ws.Range("J16").Value = 952.6
$ws.Range("K16").Value = 847.94446
$ws.Range("L16").Value = 952.6
$ws.Range("M16").Value = -560.94446
$ws.Range("N16").Value = -1526.6
# Row 31
$ws.Range("H31").Value = 3626693.2
$ws.Range("I31").Value = 2644.1072
$ws.Range("J31").Value = 9264103
$ws.Range("K31").Value = 2644.1072
$ws.Range("L31").Value = 9264103
$ws.Range("M31").Value = -2349.1072
$ws.Range("N31").Value = -9264693
# Row 34
$ws.Range("H34").Value = 3626693.2
$ws.Range("I34").Value = 2644.1072
$ws.Range("J34").Value = 9264103
$ws.Range("K34").Value = 2644.1072
$ws.Range("L34").Value = 9264103
$ws.Range("M34").Value = -2442.1072
$ws.Range("N34").Value = -9264507
# Row 58
$ws.Range("H58").Value = 1832.6666
$ws.Range("I58").Value = 2550
$ws.Range("J58").Value = 1258.8
$ws.Range("K58").Value = 2550
$ws.Range("L58").Value = 1258.8
$ws.Range("M58").Value = -2347
$ws.Range("N58").Value = -1664.8
# Row 62
$ws.Range("H62").Value = 25643128
$ws.Range("J62").Value = 66668930
$ws.Range("L62").Value = 66668930
$ws.Range("N62").Value = -66670178
# Row 65
$ws.Range("H65").Value = 25643128
$ws.Range("J65").Value = 66668930
$ws.Range("L65").Value = 333344650
$ws.Range("N65").Value = -333350890
# Row 99
$ws.Range("H99").Value = 2305.1738
$ws.Range("I99").Value = 2144.7856
$ws.Range("J99").Value = 2554.6667
$ws.Range("K99").Value = 2144.7856
$ws.Range("L99").Value = 2554.6667
$ws.Range("M99").Value = -646.7856000000002
$ws.Range("N99").Value = -5550.6667
# Row 107
$ws.Range("H107").Value = 728.60974
$ws.Range("I107").Value = 765.9048
$ws.Range("K107").Value = 765.9048
$ws.Range("M107").Value = 1154.0952
# Row 113
$ws.Range("H113").Value = 870.6957
$ws.Range("I113").Value = 847.94446
$ws.Range("J113").Value = 952.6
$ws.Range("K113").Value = 847.94446
$ws.Range("L113").Value = 952.6
$ws.Range("M113").Value = 1322.05554
$ws.Range("N113").Value = -5292.6
# Row 126
$ws.Range("H126").Value = 2305.1738
$ws.Range("I126").Value = 2144.7856
$ws.Range("J126").Value = 2554.6667
$ws.Range("K126").Value = 6434.3568
$ws.Range("L126").Value = 7664.000100000001
$ws.Range("M126").Value = -3964.3568
$ws.Range("N126").Value = -12604.0001
# Row 132
$ws.Range("H132").Value = 5956197.5
$ws.Range("I132").Value = 4554
$ws.Range("J132").Value = 8336855
$ws.Range("K132").Value = 13662
$ws.Range("L132").Value = 25010565
$ws.Range("M132").Value = -11132
$ws.Range("N132").Value = -25015625
# Row 136
$ws.Range("H136").Value = 1832.6666
$ws.Range("I136").Value = 2550
$ws.Range("J136").Value = 1258.8
$ws.Range("K136").Value = 7650
$ws.Range("L136").Value = 3776.4
$ws.Range("M136").Value = -5100
$ws.Range("N136").Value = -8876.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 761.77
$ws.Range("I131").Value = 404.44446
$ws.Range("J131").Value = 797.10986
$ws.Range("K131").Value = 1213.33338
$ws.Range("L131").Value = 2391.32958
$ws.Range("M131").Value = 3826.66662
$ws.Range("N131").Value = -12471.32958

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 5433.3335
$ws.Range("I126").Value = 6900
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 20700
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -18230
$ws.Range("N126").Value = -12440

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 33335442
$ws.Range("I7").Value = 1277.75
$ws.Range("J7").Value = 45456956
$ws.Range("K7").Value = 1277.75
$ws.Range("L7").Value = 45456956
$ws.Range("M7").Value = -1165.75
$ws.Range("N7").Value = -45457180
# Row 126
$ws.Range("H126").Value = 33335442
$ws.Range("I126").Value = 1277.75
$ws.Range("J126").Value = 45456956
$ws.Range("K126").Value = 3833.25
$ws.Range("L126").Value = 136370868
$ws.Range("M126").Value = -1363.25
$ws.Range("N126").Value = -136375808
# Row 132
$ws.Range("H132").Value = 2478.8125
$ws.Range("I132").Value = 1798
$ws.Range("J132").Value = 3159.625
$ws.Range("K132").Value = 5394
$ws.Range("L132").Value = 9478.875
$ws.Range("M132").Value = -2864
$ws.Range("N132").Value = -14538.875
# Row 136
$ws.Range("H136").Value = 9940.532999999999
$ws.Range("I136").Value = 15513.5
$ws.Range("J136").Value = 3571.4285
$ws.Range("K136").Value = 46540.5
$ws.Range("L136").Value = 10714.2855
$ws.Range("M136").Value = -43990.5
$ws.Range("N136").Value = -15814.2855

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 2897.125
$ws.Range("I136").Value = 2925.2856
$ws.Range("J136").Value = 2700
$ws.Range("K136").Value = 8775.856800000001
$ws.Range("L136").Value = 8100
$ws.Range("M136").Value = -6225.856800000001
$ws.Range("N136").Value = -13200

